$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.04806415330644143
$ws.Cells.Item(2, 4).Value = 0.09433389003166326
$ws.Cells.Item(2, 5).Value = 0.1126764088172365
$ws.Cells.Item(2, 6).Value = 2.101210578912017
$ws.Cells.Item(2, 7).Value = 1.470874964997947
$ws.Cells.Item(2, 8).Value = 1.322835871731598
$ws.Cells.Item(2, 10).Value = 0.153413385524626
$ws.Cells.Item(2, 11).Value = 1.579846439231574
$ws.Cells.Item(2, 13).Value = 0.4789933877275772
$ws.Cells.Item(2, 14).Value = 1.62875754539281

$ws.Cells.Item(3, 3).Value = 0.04273327228797541
$ws.Cells.Item(3, 4).Value = 0.09365309735122906
$ws.Cells.Item(3, 5).Value = 0.112628111743879
$ws.Cells.Item(3, 6).Value = 2.094167621395798
$ws.Cells.Item(3, 7).Value = 1.460794342100428
$ws.Cells.Item(3, 8).Value = 1.324659984804342
$ws.Cells.Item(3, 10).Value = 0.1539073061327443
$ws.Cells.Item(3, 11).Value = 1.453637886189426
$ws.Cells.Item(3, 13).Value = 0.4551671866045339
$ws.Cells.Item(3, 14).Value = 1.649415010051476

$ws.Cells.Item(4, 3).Value = 0.03947702808072506
$ws.Cells.Item(4, 4).Value = 0.09325481736216901
$ws.Cells.Item(4, 5).Value = 0.1126393877610958
$ws.Cells.Item(4, 6).Value = 2.091119198627311
$ws.Cells.Item(4, 7).Value = 1.455643872537124
$ws.Cells.Item(4, 8).Value = 1.326500421825997
$ws.Cells.Item(4, 10).Value = 0.1542801048651512
$ws.Cells.Item(4, 11).Value = 1.376727107748167
$ws.Cells.Item(4, 13).Value = 0.4407586108317645
$ws.Cells.Item(4, 14).Value = 1.662728356794245

$ws.Cells.Item(5, 3).Value = 0.0381542725570938
$ws.Cells.Item(5, 4).Value = 0.09309750690345098
$ws.Cells.Item(5, 5).Value = 0.1126542895426059
$ws.Cells.Item(5, 6).Value = 2.090197141497612
$ws.Cells.Item(5, 7).Value = 1.453805416754307
$ws.Cells.Item(5, 8).Value = 1.327431307673081
$ws.Cells.Item(5, 10).Value = 0.1544494995946444
$ws.Cells.Item(5, 11).Value = 1.345531788216135
$ws.Cells.Item(5, 13).Value = 0.4349425865951488
$ws.Cells.Item(5, 14).Value = 1.668311857775175

$ws.Cells.Item(6, 3).Value = 0.03793488110241583
$ws.Cells.Item(6, 4).Value = 0.09307168802837396
$ws.Cells.Item(6, 5).Value = 0.1126573868996612
$ws.Cells.Item(6, 6).Value = 2.09006335486724
$ws.Cells.Item(6, 7).Value = 1.453515844587812
$ws.Cells.Item(6, 8).Value = 1.327596798407697
$ws.Cells.Item(6, 10).Value = 0.1544786827684135
$ws.Cells.Item(6, 11).Value = 1.340360679276557
$ws.Cells.Item(6, 13).Value = 0.4339801995227006
$ws.Cells.Item(6, 14).Value = 1.669248545912467

$ws.Cells.Item(7, 3).Value = 0.03945917206075933
$ws.Cells.Item(7, 4).Value = 0.0932526755650116
$ws.Cells.Item(7, 5).Value = 0.1126395469808177
$ws.Cells.Item(7, 6).Value = 2.091105467842183
$ws.Cells.Item(7, 7).Value = 1.455618025291898
$ws.Cells.Item(7, 8).Value = 1.326512244007418
$ws.Cells.Item(7, 10).Value = 0.154282318631207
$ws.Cells.Item(7, 11).Value = 1.376305803585069
$ws.Cells.Item(7, 13).Value = 0.4406799487372552
$ws.Cells.Item(7, 14).Value = 1.662803017429631

$ws.Cells.Item(8, 3).Value = 0.0462225068562816
$ws.Cells.Item(8, 4).Value = 0.09409507589378308
$ws.Cells.Item(8, 5).Value = 0.1126512691655517
$ws.Cells.Item(8, 6).Value = 2.098516932392471
$ws.Cells.Item(8, 7).Value = 1.467182967244185
$ws.Cells.Item(8, 8).Value = 1.323315139131225
$ws.Cells.Item(8, 10).Value = 0.1535692536676869
$ws.Cells.Item(8, 11).Value = 1.536209026079803
$ws.Cells.Item(8, 13).Value = 0.4707323133521939
$ws.Cells.Item(8, 14).Value = 1.635749505463997

$ws.Cells.Item(9, 3).Value = 0.05962356624128518
$ws.Cells.Item(9, 4).Value = 0.09590235504150257
$ws.Cells.Item(9, 5).Value = 0.1129985049790534
$ws.Cells.Item(9, 6).Value = 2.123208664910436
$ws.Cells.Item(9, 7).Value = 1.498149617904772
$ws.Cells.Item(9, 8).Value = 1.322775341431623
$ws.Cells.Item(9, 10).Value = 0.1527230833782802
$ws.Cells.Item(9, 11).Value = 1.854406384229776
$ws.Cells.Item(9, 13).Value = 0.5314178540625178
$ws.Cells.Item(9, 14).Value = 1.587697240566941

$ws.Cells.Item(10, 3).Value = 0.06955977577099759
$ws.Cells.Item(10, 4).Value = 0.09732344005478666
$ws.Cells.Item(10, 5).Value = 0.1134507637158606
$ws.Cells.Item(10, 6).Value = 2.147594310101084
$ws.Cells.Item(10, 7).Value = 1.526018242556745
$ws.Cells.Item(10, 8).Value = 1.325892677077263
$ws.Cells.Item(10, 10).Value = 0.1524388253906324
$ws.Cells.Item(10, 11).Value = 2.091049401603698
$ws.Cells.Item(10, 13).Value = 0.5770789296284136
$ws.Cells.Item(10, 14).Value = 1.555445793669024

$ws.Cells.Item(11, 3).Value = 0.07410105570008341
$ws.Cells.Item(11, 4).Value = 0.09798989955769599
$ws.Cells.Item(11, 5).Value = 0.1136992199081845
$ws.Cells.Item(11, 6).Value = 2.16005565747615
$ws.Cells.Item(11, 7).Value = 1.539822018105411
$ws.Cells.Item(11, 8).Value = 1.328078588575949
$ws.Cells.Item(11, 10).Value = 0.1523829728534096
$ws.Cells.Item(11, 11).Value = 2.19933719099663
$ws.Cells.Item(11, 13).Value = 0.598086705113019
$ws.Cells.Item(11, 14).Value = 1.541437688631476

$ws.Cells.Item(12, 3).Value = 0.07582386788264728
$ws.Cells.Item(12, 4).Value = 0.09824511940898617
$ws.Cells.Item(12, 5).Value = 0.1137994365551869
$ws.Cells.Item(12, 6).Value = 2.164972053843798
$ws.Cells.Item(12, 7).Value = 1.545212181906095
$ws.Cells.Item(12, 8).Value = 1.329017105162109
$ws.Cells.Item(12, 10).Value = 0.1523724002204006
$ws.Cells.Item(12, 11).Value = 2.240435011418811
$ws.Cells.Item(12, 13).Value = 0.6060758278657801
$ws.Cells.Item(12, 14).Value = 1.536228702149288

$ws.Cells.Item(13, 3).Value = 0.07545268917372994
$ws.Cells.Item(13, 4).Value = 0.09819002714074543
$ws.Cells.Item(13, 5).Value = 0.1137775805814201
$ws.Cells.Item(13, 6).Value = 2.163904420628228
$ws.Cells.Item(13, 7).Value = 1.544044048936456
$ws.Cells.Item(13, 8).Value = 1.328810047175551
$ws.Cells.Item(13, 10).Value = 0.1523742065553435
$ws.Cells.Item(13, 11).Value = 2.231579788817044
$ws.Cells.Item(13, 13).Value = 0.604353718840855
$ws.Cells.Item(13, 14).Value = 1.537346296336642

$ws.Cells.Item(14, 3).Value = 0.07424272938717991
$ws.Cells.Item(14, 4).Value = 0.09801083977832548
$ws.Cells.Item(14, 5).Value = 0.1137073419696257
$ws.Cells.Item(14, 6).Value = 2.160456168404764
$ws.Cells.Item(14, 7).Value = 1.540262197975039
$ws.Cells.Item(14, 8).Value = 1.328153579107692
$ws.Cells.Item(14, 10).Value = 0.1523818909635963
$ws.Cells.Item(14, 11).Value = 2.202716498510085
$ws.Cells.Item(14, 13).Value = 0.5987432952247218
$ws.Cells.Item(14, 14).Value = 1.541007225432168

$ws.Cells.Item(15, 3).Value = 0.07350200360548342
$ws.Cells.Item(15, 4).Value = 0.09790145210804013
$ws.Cells.Item(15, 5).Value = 0.1136651169771525
$ws.Cells.Item(15, 6).Value = 2.15836976598213
$ws.Cells.Item(15, 7).Value = 1.537966959680404
$ws.Cells.Item(15, 8).Value = 1.327765908005802
$ws.Cells.Item(15, 10).Value = 0.1523879758229256
$ws.Cells.Item(15, 11).Value = 2.185048823751629
$ws.Cells.Item(15, 13).Value = 0.5953111666072886
$ws.Cells.Item(15, 14).Value = 1.543262103629957

$ws.Cells.Item(16, 3).Value = 0.06926342620830894
$ws.Cells.Item(16, 4).Value = 0.09728028525446319
$ws.Cells.Item(16, 5).Value = 0.1134353851855394
$ws.Cells.Item(16, 6).Value = 2.146807532331081
$ws.Cells.Item(16, 7).Value = 1.525138882997908
$ws.Cells.Item(16, 8).Value = 1.325765302729565
$ws.Cells.Item(16, 10).Value = 0.1524439548633012
$ws.Cells.Item(16, 11).Value = 2.083985361658165
$ws.Cells.Item(16, 13).Value = 0.5757107732449072
$ws.Cells.Item(16, 14).Value = 1.556374619786373

$ws.Cells.Item(17, 3).Value = 0.06666868993180231
$ws.Cells.Item(17, 4).Value = 0.09690432074921063
$ws.Cells.Item(17, 5).Value = 0.1133053858151385
$ws.Cells.Item(17, 6).Value = 2.140065477214151
$ws.Cells.Item(17, 7).Value = 1.517558448299127
$ws.Cells.Item(17, 8).Value = 1.324734899082983
$ws.Cells.Item(17, 10).Value = 0.1524971205981558
$ws.Cells.Item(17, 11).Value = 2.022149441838451
$ws.Cells.Item(17, 13).Value = 0.5637470586723623
$ws.Cells.Item(17, 14).Value = 1.564588739511251

$ws.Cells.Item(18, 3).Value = 0.06517826188411391
$ws.Cells.Item(18, 4).Value = 0.0966899602509983
$ws.Cells.Item(18, 5).Value = 0.1132346354380793
$ws.Cells.Item(18, 6).Value = 2.13631635746691
$ws.Cells.Item(18, 7).Value = 1.51330432394272
$ws.Cells.Item(18, 8).Value = 1.324214487910552
$ws.Cells.Item(18, 10).Value = 0.1525346134409347
$ws.Cells.Item(18, 11).Value = 1.986643084658056
$ws.Cells.Item(18, 13).Value = 0.5568880929929065
$ws.Cells.Item(18, 14).Value = 1.569375698466796

$ws.Cells.Item(19, 3).Value = 0.06467396910525736
$ws.Cells.Item(19, 4).Value = 0.09661770602663466
$ws.Cells.Item(19, 5).Value = 0.1132113716854164
$ws.Cells.Item(19, 6).Value = 2.135069055256309
$ws.Cells.Item(19, 7).Value = 1.511882114851147
$ws.Cells.Item(19, 8).Value = 1.32405068407121
$ws.Cells.Item(19, 10).Value = 0.1525484948114837
$ws.Cells.Item(19, 11).Value = 1.974631555959888
$ws.Cells.Item(19, 13).Value = 0.5545695874862773
$ws.Cells.Item(19, 14).Value = 1.571007194113237

$ws.Cells.Item(20, 3).Value = 0.06694469697237082
$ws.Cells.Item(20, 4).Value = 0.09694414800787143
$ws.Cells.Item(20, 5).Value = 0.1133188083118348
$ws.Cells.Item(20, 6).Value = 2.140769851413438
$ws.Cells.Item(20, 7).Value = 1.518354426612092
$ws.Cells.Item(20, 8).Value = 1.324837107258901
$ws.Cells.Item(20, 10).Value = 0.1524907454108728
$ws.Cells.Item(20, 11).Value = 2.028725769246194
$ws.Cells.Item(20, 13).Value = 0.5650183144514358
$ws.Cells.Item(20, 14).Value = 1.563707871333765

$ws.Cells.Item(21, 3).Value = 0.07459803843167379
$ws.Cells.Item(21, 4).Value = 0.09806339444546808
$ws.Cells.Item(21, 5).Value = 0.1137278064318714
$ws.Cells.Item(21, 6).Value = 2.161463635593819
$ws.Cells.Item(21, 7).Value = 1.541368587986199
$ws.Cells.Item(21, 8).Value = 1.328343390917468
$ws.Cells.Item(21, 10).Value = 0.152379346679588
$ws.Cells.Item(21, 11).Value = 2.211191861763325
$ws.Cells.Item(21, 13).Value = 0.6003902920457449
$ws.Cells.Item(21, 14).Value = 1.539929326154759

$ws.Cells.Item(22, 3).Value = 0.07961820610832149
$ws.Cells.Item(22, 4).Value = 0.09881146227635185
$ws.Cells.Item(22, 5).Value = 0.1140308452546144
$ws.Cells.Item(22, 6).Value = 2.176140101660934
$ws.Cells.Item(22, 7).Value = 1.55736009077566
$ws.Cells.Item(22, 8).Value = 1.331280705310377
$ws.Cells.Item(22, 10).Value = 0.1523682019531236
$ws.Cells.Item(22, 11).Value = 2.330977788001178
$ws.Cells.Item(22, 13).Value = 0.6237056985216611
$ws.Cells.Item(22, 14).Value = 1.524945879738153

$ws.Cells.Item(23, 3).Value = 0.07693715263498291
$ws.Cells.Item(23, 4).Value = 0.09841069754912013
$ws.Cells.Item(23, 5).Value = 0.1138658417562795
$ws.Cells.Item(23, 6).Value = 2.168201328207203
$ws.Cells.Item(23, 7).Value = 1.54873781966171
$ws.Cells.Item(23, 8).Value = 1.329653802435047
$ws.Cells.Item(23, 10).Value = 0.1523685032861835
$ws.Cells.Item(23, 11).Value = 2.26699697605369
$ws.Cells.Item(23, 13).Value = 0.6112437459506168
$ws.Cells.Item(23, 14).Value = 1.532891772722603

$ws.Cells.Item(24, 3).Value = 0.06681991008349542
$ws.Cells.Item(24, 4).Value = 0.0969261365376326
$ws.Cells.Item(24, 5).Value = 0.1133127275785419
$ws.Cells.Item(24, 6).Value = 2.140451008375194
$ws.Cells.Item(24, 7).Value = 1.51799424104351
$ws.Cells.Item(24, 8).Value = 1.324790674749323
$ws.Cells.Item(24, 10).Value = 0.1524936060574404
$ws.Cells.Item(24, 11).Value = 2.025752474762839
$ws.Cells.Item(24, 13).Value = 0.5644435201136915
$ws.Cells.Item(24, 14).Value = 1.564105910773982

$ws.Cells.Item(25, 3).Value = 0.05598279107107373
$ws.Cells.Item(25, 4).Value = 0.09539695647535495
$ws.Cells.Item(25, 5).Value = 0.1128699107286586
$ws.Cells.Item(25, 6).Value = 2.115436073399167
$ws.Cells.Item(25, 7).Value = 1.488878366640876
$ws.Cells.Item(25, 8).Value = 1.322305819887788
$ws.Cells.Item(25, 10).Value = 0.1528928004103385
$ws.Cells.Item(25, 11).Value = 1.767825848585517
$ws.Cells.Item(25, 13).Value = 0.5148124587319458
$ws.Cells.Item(25, 14).Value = 1.600160955693608

